$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-18 Friday" "2024-10-19 Saturday"
Replace-Text "886×4=" "540×9="
Replace-Text "595×8=" "547×9="
Replace-Text "448×7=" "890×7="
Replace-Text "361×6=" "791×5="
Replace-Text "980×5=" "632×8="
Replace-Text "175×5=" "838×9="
Replace-Text "304×4=" "815×5="
Replace-Text "551×9=" "253×2="
Replace-Text "705×9=" "615×4="
Replace-Text "551×7=" "219×3="
Replace-Text "718×6=" "844×9="
Replace-Text "120×3=" "437×9="
Replace-Text "405×2=" "356×8="
Replace-Text "286×7=" "447×8="
Replace-Text "692×2=" "609×7="
Replace-Text "214×5=" "710×5="
Replace-Text "776×4=" "257×4="
Replace-Text "504×6=" "204×3="
Replace-Text "597×3=" "110×9="
Replace-Text "955×3=" "546×3="
Replace-Text "916×6=" "508×5="
Replace-Text "775×2=" "928×7="
Replace-Text "477×3=" "837×3="
Replace-Text "916×2=" "716×7="
Replace-Text "225×3=" "848×3="
